$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -11.7662
$ws.Range("B3").Value = 5.221500000000003
$ws.Range("C3").Value = -13.1476
$ws.Range("E3").Value = 16.56030000000001
$ws.Range("B4").Value = 9.212600000000002
$ws.Range("E5").Value = 16.5782
$ws.Range("E6").Value = 16.3122
$ws.Range("D8").Value = -8.528600000000004
$ws.Range("C9").Value = -10.2936
$ws.Range("E10").Value = 16.6567
$ws.Range("A11").Value = -21.7927
$ws.Range("D11").Value = -7.140799999999996
$ws.Range("A12").Value = -21.62830000000001
$ws.Range("B14").Value = 6.935499999999995
$ws.Range("D14").Value = -7.4079
$ws.Range("A15").Value = -21.87289999999998
$ws.Range("C15").Value = -13.21029999999999
$ws.Range("D15").Value = -8.771399999999995
$ws.Range("D17").Value = -8.1883
$ws.Range("C19").Value = -12.4396
$ws.Range("C20").Value = -12.659
$ws.Range("E21").Value = 16.56730000000001
$ws.Range("C25").Value = -13.14040000000001
$ws.Range("B26").Value = 4.794900000000001
$ws.Range("D26").Value = -8.441400000000005
$ws.Range("A27").Value = -21.63659999999998
$ws.Range("C27").Value = -13.64329999999999
$ws.Range("E27").Value = 16.4905
$ws.Range("A28").Value = -21.70979999999999
$ws.Range("C28").Value = -13.2981
$ws.Range("E29").Value = 16.98880000000001
$ws.Range("C30").Value = -13.37209999999999
$ws.Range("A31").Value = -21.2886
$ws.Range("B31").Value = 4.690100000000002
$ws.Range("A32").Value = -21.70020000000001
$ws.Range("C32").Value = -13.072
$ws.Range("E33").Value = 17.04680000000002
$ws.Range("B35").Value = 8.653400000000003
$ws.Range("A36").Value = -20.1192
$ws.Range("D36").Value = -7.151500000000002
$ws.Range("E36").Value = 16.98950000000001
$ws.Range("B37").Value = 8.560400000000001
$ws.Range("A38").Value = -18.93249999999999
$ws.Range("B39").Value = 8.884699999999999
$ws.Range("E39").Value = 15.77769999999999
$ws.Range("B40").Value = 8.746599999999995
$ws.Range("D42").Value = -9.007299999999997
$ws.Range("C44").Value = -14.00499999999999
$ws.Range("B45").Value = 6.720899999999995
$ws.Range("A46").Value = -21.81100000000002
$ws.Range("C47").Value = -12.3457
$ws.Range("E47").Value = 16.7022
$ws.Range("B52").Value = 5.233400000000001
$ws.Range("E53").Value = 16.4403
$ws.Range("A54").Value = -22.11879999999999
$ws.Range("E54").Value = 16.7916
$ws.Range("A55").Value = -22.27529999999999
$ws.Range("A56").Value = -21.71090000000001
$ws.Range("E56").Value = 15.9272
$ws.Range("B57").Value = 4.731799999999996
$ws.Range("C58").Value = -13.758
$ws.Range("E58").Value = 16.1675
$ws.Range("E60").Value = 15.39770000000001
$ws.Range("C62").Value = -14.4793
$ws.Range("D64").Value = -7.238999999999994
$ws.Range("E66").Value = 17.00580000000002
$ws.Range("A67").Value = -21.56539999999998
$ws.Range("D68").Value = -7.006499999999993
$ws.Range("A69").Value = -21.78809999999999
$ws.Range("E69").Value = 17.35310000000001
$ws.Range("A72").Value = -22.11650000000002
$ws.Range("E72").Value = 17.05019999999999
$ws.Range("A73").Value = -20.11789999999999
$ws.Range("C77").Value = -11.08159999999999
$ws.Range("C78").Value = -10.938
$ws.Range("D79").Value = -5.724499999999997
$ws.Range("E80").Value = 16.7871
$ws.Range("B81").Value = 6.557199999999999
$ws.Range("E82").Value = 16.77150000000001
$ws.Range("A83").Value = -21.88039999999999
$ws.Range("B83").Value = 5.919800000000001
$ws.Range("E83").Value = 16.70470000000001
$ws.Range("C84").Value = -13.79139999999999
$ws.Range("A86").Value = -22.1541
$ws.Range("C89").Value = -10.27639999999999
$ws.Range("D89").Value = -5.626099999999999
$ws.Range("A91").Value = -21.41960000000001
$ws.Range("C91").Value = -10.38729999999999
$ws.Range("C92").Value = -10.80399999999999
$ws.Range("A93").Value = -21.3119
$ws.Range("C96").Value = -13.6806
$ws.Range("A99").Value = -20.34909999999999
$ws.Range("B100").Value = 4.389699999999999
$ws.Range("B102").Value = 8.104399999999998
$ws.Range("C102").Value = -13.5349
